$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the floating profile-photo paragraph (image2.jpeg anchored
#    drawing) together with the run of empty spacer paragraphs that follow
#    it, right before the "BRHANE GIDAY" Title paragraph.
# ---------------------------------------------------------------------------

# Locate the Title paragraph ("BRHANE GIDAY") dynamically so the script does
# not depend on a brittle hard-coded paragraph index.
$titleIndex = -1
$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($p.Style.NameLocal -eq "Title") {
        $titleIndex = $idx
        break
    }
}

if ($titleIndex -gt 1) {
    $firstToRemove = $d.Paragraphs.Item(1)
    $lastToRemove = $d.Paragraphs.Item($titleIndex - 1)
    # Walk backwards from the Title paragraph to find the first paragraph of
    # the contiguous block of "BodyText"-styled spacer paragraphs (the
    # photo paragraph + the blank paragraphs under it) that sits right above
    # the Title paragraph.
    $blockStart = $titleIndex - 1
    while ($blockStart -gt 1) {
        $prev = $d.Paragraphs.Item($blockStart - 1)
        if ($prev.Style.NameLocal -ne "Body Text") {
            break
        }
        $blockStart = $blockStart - 1
    }
    $startPara = $d.Paragraphs.Item($blockStart)
    $endPara = $d.Paragraphs.Item($titleIndex - 1)
    $rng = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $rng.Delete()
}

# ---------------------------------------------------------------------------
# 2) Remove the "High school: Axum Secondary school (2010-2014)" bullet
#    paragraph from the EDUCATION section entirely.
# ---------------------------------------------------------------------------

$hsIndex = -1
$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($p.Range.Text -like "*High school*Axum*") {
        $hsIndex = $idx
        break
    }
}

if ($hsIndex -gt 0) {
    $hsPara = $d.Paragraphs.Item($hsIndex)
    $hsPara.Range.Delete()
}
